$wb = $excel.ActiveWorkbook

# --- paper sheet: per-country specific consumption fix ---
$ws = $wb.Worksheets.Item("paper")
$ws.Cells.Item(2, 2).ClearContents()
$ws.Cells.Item(2, 3).Value = 0.8700847
$ws.Cells.Item(2, 4).ClearContents()
$ws.Cells.Item(3, 2).ClearContents()
$ws.Cells.Item(3, 3).Value = -5.0731063
$ws.Cells.Item(3, 4).ClearContents()
$ws.Cells.Item(4, 2).ClearContents()
$ws.Cells.Item(4, 3).Value = 1.7631324
$ws.Cells.Item(4, 4).ClearContents()
$ws.Cells.Item(5, 2).Value = 2.9946903
$ws.Cells.Item(5, 3).Value = 7.4707965
$ws.Cells.Item(5, 4).Value = 0
$ws.Cells.Item(6, 2).ClearContents()
$ws.Cells.Item(6, 3).Value = 1.2620384
$ws.Cells.Item(6, 4).ClearContents()
$ws.Cells.Item(7, 2).Value = 14.8435663
$ws.Cells.Item(7, 3).ClearContents()
$ws.Cells.Item(7, 4).ClearContents()
$ws.Cells.Item(8, 2).Value = 3.0854156
$ws.Cells.Item(8, 3).ClearContents()
$ws.Cells.Item(8, 4).ClearContents()
$ws.Cells.Item(9, 2).Value = 2.0036631
$ws.Cells.Item(9, 3).ClearContents()
$ws.Cells.Item(9, 4).ClearContents()
$ws.Cells.Item(10, 2).Value = 3.8029549
$ws.Cells.Item(10, 3).ClearContents()
$ws.Cells.Item(10, 4).ClearContents()
$ws.Cells.Item(11, 2).ClearContents()
$ws.Cells.Item(11, 3).Value = -0.2147037
$ws.Cells.Item(11, 4).ClearContents()
$ws.Cells.Item(12, 2).ClearContents()
$ws.Cells.Item(12, 3).Value = 4.2137493
$ws.Cells.Item(12, 4).ClearContents()
$ws.Cells.Item(13, 2).ClearContents()
$ws.Cells.Item(13, 3).Value = -0.1167503
$ws.Cells.Item(13, 4).ClearContents()
$ws.Cells.Item(14, 2).ClearContents()
$ws.Cells.Item(14, 3).Value = 0.0916606
$ws.Cells.Item(14, 4).ClearContents()
$ws.Cells.Item(15, 2).ClearContents()
$ws.Cells.Item(15, 3).Value = 2.3427368
$ws.Cells.Item(15, 4).ClearContents()
$ws.Cells.Item(16, 2).ClearContents()
$ws.Cells.Item(16, 3).Value = 1.9525423
$ws.Cells.Item(16, 4).ClearContents()
$ws.Cells.Item(17, 2).Value = 2.9946903
$ws.Cells.Item(17, 3).Value = 7.4707965
$ws.Cells.Item(17, 4).Value = 0
$ws.Cells.Item(18, 2).ClearContents()
$ws.Cells.Item(18, 3).Value = -5.249878
$ws.Cells.Item(18, 4).ClearContents()
$ws.Cells.Item(19, 2).ClearContents()
$ws.Cells.Item(19, 3).Value = 2.4530077
$ws.Cells.Item(19, 4).ClearContents()
$ws.Cells.Item(20, 2).Value = 2.9946903
$ws.Cells.Item(20, 3).Value = 7.4707965
$ws.Cells.Item(20, 4).Value = 0
$ws.Cells.Item(21, 2).ClearContents()
$ws.Cells.Item(21, 3).Value = 0.0659781
$ws.Cells.Item(21, 4).ClearContents()
$ws.Cells.Item(22, 2).ClearContents()
$ws.Cells.Item(22, 3).Value = 2.6669077
$ws.Cells.Item(22, 4).ClearContents()
$ws.Cells.Item(23, 2).ClearContents()
$ws.Cells.Item(23, 3).Value = 1.8430915
$ws.Cells.Item(23, 4).ClearContents()
$ws.Cells.Item(24, 2).Value = 6.2365862
$ws.Cells.Item(24, 3).ClearContents()
$ws.Cells.Item(24, 4).ClearContents()
$ws.Cells.Item(25, 2).Value = 2.9946903
$ws.Cells.Item(25, 3).Value = 7.4707965
$ws.Cells.Item(25, 4).Value = 0
$ws.Cells.Item(26, 2).ClearContents()
$ws.Cells.Item(26, 3).Value = 0.0270267
$ws.Cells.Item(26, 4).ClearContents()
$ws.Cells.Item(27, 2).Value = 2.9946903
$ws.Cells.Item(27, 3).Value = 7.4707965
$ws.Cells.Item(27, 4).Value = 0
$ws.Cells.Item(28, 2).Value = 2.9946903
$ws.Cells.Item(28, 3).Value = 7.4707965
$ws.Cells.Item(28, 4).Value = 0
$ws.Cells.Item(29, 2).ClearContents()
$ws.Cells.Item(29, 3).Value = -0.2755971
$ws.Cells.Item(29, 4).ClearContents()
$ws.Cells.Item(30, 2).Value = 2.9946903
$ws.Cells.Item(30, 3).Value = 7.4707965
$ws.Cells.Item(30, 4).Value = 0
$ws.Cells.Item(31, 2).ClearContents()
$ws.Cells.Item(31, 3).Value = 1.8578062
$ws.Cells.Item(31, 4).ClearContents()
$ws.Cells.Item(32, 2).Value = 2.9946903
$ws.Cells.Item(32, 3).Value = 7.4707965
$ws.Cells.Item(32, 4).Value = 0
$ws.Cells.Item(33, 2).Value = -13.9284813
$ws.Cells.Item(33, 3).ClearContents()
$ws.Cells.Item(33, 4).ClearContents()
$ws.Cells.Item(34, 2).ClearContents()
$ws.Cells.Item(34, 3).Value = -0.7921223
$ws.Cells.Item(34, 4).ClearContents()
$ws.Cells.Item(35, 2).ClearContents()
$ws.Cells.Item(35, 3).Value = 0.4206817
$ws.Cells.Item(35, 4).ClearContents()

# --- cement sheet: normalize Electricity/Heat values, reset H2 substitution cap to 0 ---
$ws = $wb.Worksheets.Item("cement")
$ws.Range("B2:B35").Value = 0.42092
$ws.Range("C2:C35").Value = 2.97908
$ws.Range("E2:E35").Value = 0

# --- glass sheet: normalize Electricity/Heat values, reset H2 substitution cap to 0 ---
$ws = $wb.Worksheets.Item("glass")
$ws.Range("B2:B35").Value = 1.3898961
$ws.Range("C2:C35").Value = 6.1331626
$ws.Range("E2:E35").Value = 0

# --- alu_prim, alu_sec, copper_prim, copper_sec: reset H2 substitution cap to 0 ---
$ws = $wb.Worksheets.Item("alu_prim")
$ws.Range("E2:E35").Value = 0
$ws = $wb.Worksheets.Item("alu_sec")
$ws.Range("E2:E35").Value = 0
$ws = $wb.Worksheets.Item("copper_prim")
$ws.Range("E2:E35").Value = 0
$ws = $wb.Worksheets.Item("copper_sec")
$ws.Range("E2:E35").Value = 0
